# update version and docs
# - add 4 new function-translation rows/cells to the "export" sheet
# - update the saved view (scroll position + selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# square(x): SBML translation changes from "x^2" to "power(x,2)"
$ws.Range("E19").Value = "power(x,2)"

# New DBSolve-alternative rows gain a Heta-column (B) label
$ws.Range("B30").Value = "ifg0(x-y, 1, 2)"
$ws.Range("B31").Value = "ifge0(x-y, 1, 2)"
$ws.Range("B34").Value = "ife0(x-y, 1, 2)"

# Restore the view: scrolled down so row 19 is at the top, with F34 selected
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F34").Select()
